$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep numeric-looking Price values as text (matches source formatting)
$textCells = 'D5','D6','D11','D14','D18','D19','D21','D22','D23','D24','D26','D28','D31','D32','D34','D35','D36','D37','D41','D44','D46','D47','D48','D49','D50','D51'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data
$ws.Range('D2').Value = '70.954.84'
$ws.Range('E2').Value = '  -0.68%  '

$ws.Range('D3').Value = '3.796.29'
$ws.Range('E3').Value = '  -1.19%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '699.11'
$ws.Range('E5').Value = '  -1.32%  '

$ws.Range('D6').Value = '169.51'
$ws.Range('E6').Value = '  -2.20%  '

$ws.Range('D7').Value = '3.798.21'
$ws.Range('E7').Value = '  -1.14%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  -0.54%  '

$ws.Range('D11').Value = '7.52'
$ws.Range('E11').Value = '  +2.56%  '

$ws.Range('E12').Value = '  +4.21%  '

$ws.Range('D14').Value = '36.18'
$ws.Range('E14').Value = '  -2.36%  '

$ws.Range('D15').Value = '4.441.82'
$ws.Range('E15').Value = '  -1.07%  '

$ws.Range('D16').Value = '3.812.42'
$ws.Range('E16').Value = '  -0.87%  '

$ws.Range('D17').Value = '71.124.20'
$ws.Range('E17').Value = '  -0.47%  '

$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '7.21'
$ws.Range('E18').Value = '  -0.92%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '17.59'
$ws.Range('E19').Value = '  +0.47%  '

$ws.Range('E20').Value = '  +0.17%  '

$ws.Range('D21').Value = '512.06'
$ws.Range('E21').Value = '  +2.67%  '

$ws.Range('D22').Value = '10.44'
$ws.Range('E22').Value = '  -2.74%  '

$ws.Range('D23').Value = '0.715'
$ws.Range('E23').Value = '  -2.65%  '

$ws.Range('D24').Value = '83.53'
$ws.Range('E24').Value = '  -2.36%  '

$ws.Range('E25').Value = '  -3.62%  '

$ws.Range('D26').Value = '12.64'
$ws.Range('E26').Value = '  +3.14%  '

$ws.Range('D27').Value = '3.947.40'
$ws.Range('E27').Value = '  -1.18%  '

$ws.Range('D28').Value = '10.25'
$ws.Range('E28').Value = '  -4.60%  '

$ws.Range('E29').Value = '  -0.01%  '

$ws.Range('E30').Value = '  -5.61%  '

$ws.Range('D31').Value = '2.98'
$ws.Range('E31').Value = '  -4.76%  '

$ws.Range('D32').Value = '2.26'
$ws.Range('E32').Value = '  +0.69%  '

$ws.Range('E33').Value = '  -2.98%  '

$ws.Range('D34').Value = '29.07'
$ws.Range('E34').Value = '  -1.46%  '

$ws.Range('D35').Value = '0.171'
$ws.Range('E35').Value = '  -4.96%  '

$ws.Range('D36').Value = '9.30'
$ws.Range('E36').Value = '  +0.51%  '

$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.14%  '

$ws.Range('B38').Value = 'RenzoRestakedETH'
$ws.Range('C38').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D38').Value = '3.765.95'
$ws.Range('E38').Value = '  -1.05%  '

$ws.Range('E39').Value = '  +10.02%  '

$ws.Range('E40').Value = '  -2.84%  '

$ws.Range('D41').Value = '2.36'
$ws.Range('E41').Value = '  -0.08%  '

$ws.Range('D44').Value = '3.17'
$ws.Range('E44').Value = '  -6.34%  '

$ws.Range('E45').Value = '  -0.08%  '

$ws.Range('D46').Value = '164.02'
$ws.Range('E46').Value = '  -0.07%  '

$ws.Range('D47').Value = '49.28'
$ws.Range('E47').Value = '  +0.49%  '

$ws.Range('D48').Value = '0.000304'
$ws.Range('E48').Value = '  -4.43%  '

$ws.Range('D49').Value = '421.87'
$ws.Range('E49').Value = '  -2.21%  '

$ws.Range('D50').Value = '8.62'
$ws.Range('E50').Value = '  -1.60%  '

$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').Value = '1.36'
$ws.Range('E51').Value = '  -1.68%  '
